$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table 1 additions: A18:B22 (Num Angle / Fluence) ---
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = 1260000000000000
$ws.Range("A19").Value = 8
$ws.Range("B19").Value = 1360000000000000
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 1390000000000000
$ws.Range("A21").Value = 14
$ws.Range("B21").Value = 1380000000000000
$ws.Range("A22").Value = 17
$ws.Range("B22").Value = 1380000000000000000

$ws.Range("B18:B22").NumberFormat = "0.00E+00"

# --- Table 2 additions: A24:C31 (Num Angle / Fluence / Normalized) ---
$ws.Range("A24").Value = 25
$ws.Range("B24").Value = 9820000000000000
$ws.Range("A25").Value = 50
$ws.Range("B25").Value = 20300000000000000
$ws.Range("A26").Value = 75
$ws.Range("B26").Value = 30400000000000000
$ws.Range("A27").Value = 100
$ws.Range("B27").Value = 41000000000000000
$ws.Range("A28").Value = 125
$ws.Range("B28").Value = 51100000000000000
$ws.Range("A29").Value = 150
$ws.Range("B29").Value = 61400000000000000
$ws.Range("A30").Value = 175
$ws.Range("B30").Value = 71700000000000000
$ws.Range("A31").Value = 200
$ws.Range("B31").Value = 81700000000000000

$ws.Range("C24").Formula = "=B24/B`$24"
$ws.Range("C25:C31").Formula = "=B25/B`$24"

$ws.Range("B24:B31").NumberFormat = "0.00E+00"
$ws.Range("C24:C31").NumberFormat = "0.00E+00"

# --- View state: scroll position + selection ---
$ws.Range("B32").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
